$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-interpreted by Excel as numbers (these were/are stored as text).
$textForcedCells = @("D4", "D5", "D6", "D9", "D10", "D13", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.807.64'
$ws.Range("E2").Value = '  -4.13%  '
$ws.Range("D3").Value = '3.650.87'
$ws.Range("E3").Value = '  -3.90%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '590.64'
$ws.Range("E5").Value = '  -3.75%  '
$ws.Range("D6").Value = '163.62'
$ws.Range("E6").Value = '  -7.61%  '
$ws.Range("D7").Value = '3.647.37'
$ws.Range("E7").Value = '  -4.01%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").Value = '0.157'
$ws.Range("E10").Value = '  -5.77%  '
$ws.Range("E11").Value = '  -6.29%  '
$ws.Range("E12").Value = '  -5.55%  '
$ws.Range("D13").Value = '37.10'
$ws.Range("E13").Value = '  -6.74%  '
$ws.Range("E14").Value = '  -6.93%  '
$ws.Range("D15").Value = '4.252.96'
$ws.Range("E15").Value = '  -4.15%  '
$ws.Range("D16").Value = '3.653.75'
$ws.Range("E16").Value = '  -4.03%  '
$ws.Range("D17").Value = '66.810.22'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("D18").Value = '0.113'
$ws.Range("E18").Value = '  -4.52%  '
$ws.Range("D19").Value = '7.07'
$ws.Range("E19").Value = '  -6.38%  '
$ws.Range("D20").Value = '16.82'
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '485.72'
$ws.Range("E21").Value = '  -4.34%  '
$ws.Range("D22").Value = '8.97'
$ws.Range("E22").Value = '  -7.09%  '
$ws.Range("D23").Value = '0.708'
$ws.Range("E23").Value = '  -3.86%  '
$ws.Range("E24").Value = '  -1.74%  '
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -8.30%  '
$ws.Range("D26").Value = '0.0000136'
$ws.Range("E26").Value = '  -5.62%  '
$ws.Range("D27").Value = '12.02'
$ws.Range("E27").Value = '  -5.38%  '
$ws.Range("D28").Value = '0.993'
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("D29").Value = '9.85'
$ws.Range("E29").Value = '  -6.62%  '
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("E31").Value = '  -7.34%  '
$ws.Range("D32").Value = '7.63'
$ws.Range("E32").Value = '  -5.41%  '
$ws.Range("D33").Value = '31.34'
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").Value = '3.784.03'
$ws.Range("E34").Value = '  -4.13%  '
$ws.Range("D35").Value = '3.579.96'
$ws.Range("E35").Value = '  -4.19%  '
$ws.Range("E36").Value = '  -7.39%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '0.986'
$ws.Range("E38").Value = '  -5.65%  '
$ws.Range("D39").Value = '5.69'
$ws.Range("E39").Value = '  -6.96%  '
$ws.Range("E40").Value = '  -8.13%  '
$ws.Range("D41").Value = '0.319'
$ws.Range("E41").Value = '  -5.81%  '
$ws.Range("D42").Value = '430.64'
$ws.Range("E42").Value = '  -10.48%  '
$ws.Range("D43").Value = '48.39'
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("D44").Value = '1.90'
$ws.Range("E44").Value = '  -7.82%  '
$ws.Range("D45").Value = '2.75'
$ws.Range("E45").Value = '  -8.96%  '
$ws.Range("D46").Value = '8.26'
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("D48").Value = '142.09'
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("D49").Value = '39.46'
$ws.Range("E49").Value = '  -10.68%  '
$ws.Range("D50").Value = '2.737.06'
$ws.Range("E50").Value = '  -6.65%  '
$ws.Range("D51").Value = '0.0342'
$ws.Range("E51").Value = '  -5.77%  '
